$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-24
$gValues = @{
    2  = 1
    3  = 3
    4  = 9
    5  = 2
    6  = 6
    7  = 5
    8  = 5
    9  = 5
    10 = 0
    11 = 1
    12 = 3
    13 = 3
    14 = 1
    15 = 2
    16 = 0
    17 = 2
    18 = 4
    19 = 3
    20 = 5
    21 = 2
    22 = 1
    23 = 3
    24 = 3
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
